$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "U" with "T" in every codon string found in column A (rows 2-65),
# mirroring an RNA->DNA notation fix (Find & Replace on column A).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val.ToString().Replace("U", "T")
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Reflect the resulting selection left behind by the find & replace across column A.
$ws.Range("A1:A1048576").Select()
